$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 152, shifting rows 152:227 down to 153:228
$ws.Rows.Item(152).Insert()

# Populate the newly inserted row 152 with the new data record
$ws.Cells.Item(152, 1).Value = 3
$ws.Cells.Item(152, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(152, 3).Value = "Coquimbo"
$ws.Cells.Item(152, 4).Value = 44839
$ws.Cells.Item(152, 5).Value = 5
$ws.Cells.Item(152, 6).Value = 100112010
$ws.Cells.Item(152, 7).Value = "Achicoria"
$ws.Cells.Item(152, 8).Value = "Sin especificar"
$ws.Cells.Item(152, 9).Value = "Primera"
$ws.Cells.Item(152, 10).Value = 60
$ws.Cells.Item(152, 11).Value = 6000
$ws.Cells.Item(152, 12).Value = 6000
$ws.Cells.Item(152, 13).Value = 6000
$ws.Cells.Item(152, 14).Value = "`$/caja 16 unidades"
$ws.Cells.Item(152, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(152, 16).Value = 375
$ws.Cells.Item(152, 17).Value = 16
$ws.Cells.Item(152, 18).Value = "Hortaliza"
